$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# C2/D2 previously held numeric "Lead ID" values (10047/10049); update them to
# text values "leaf"/"Indhu", entered with a leading apostrophe so Excel keeps
# them quote-prefixed (text-forced) while keeping the center-aligned look.
$ws.Range("C2").Value = "'leaf"
$ws.Range("D2").Value = "'Indhu"

# Move the active selection to D2
$ws.Range("D2").Select()
